$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Förändrad" date (column C) bumped by one day for every data row (2-11): 46073 -> 46074
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 46074
}

# 2) Rows 7-11 got re-ordered (the logging-notice rows rotated up by one; the
#    row that used to be first (row 7) wrapped around to the bottom (row 11)).
#    Capture the "before" values of columns A (Beteckning), B (Datum) and
#    G (Area (ha)) first, then write them back shifted.
$colA = @{}
$colB = @{}
$colG = @{}
for ($r = 7; $r -le 11; $r++) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value2
    $colB[$r] = $ws.Cells.Item($r, 2).Value2
    $colG[$r] = $ws.Cells.Item($r, 7).Value2
}

for ($r = 7; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $colA[$r + 1]
    $ws.Cells.Item($r, 2).Value2 = $colB[$r + 1]
    $ws.Cells.Item($r, 7).Value2 = $colG[$r + 1]
}
$ws.Cells.Item(11, 1).Value2 = $colA[7]
$ws.Cells.Item(11, 2).Value2 = $colB[7]
$ws.Cells.Item(11, 7).Value2 = $colG[7]
